$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Ingestion sheet: add four new columnstore result rows (DW9..DW12)
# ---------------------------------------------------------------------------
$ing = $wb.Worksheets.Item("Ingestion")

# Row 10 - DW9
$ing.Range("A10").Value = "DW9"
$ing.Range("B10").Value = 11794.974
$ing.Range("C10").Value = 14400
$ing.Range("C10").Font.Bold = $true
$ing.Range("C10").Font.Italic = $true
$ing.Range("D10").Value = 152140.85200000001
$ing.Range("D10").NumberFormat = "#,##0.000"
$ing.Range("E10").Value = 10320
$ing.Range("F10").Value = 59608.641000000003
$ing.Range("F10").NumberFormat = "#,##0.00"
$ing.Range("G10").Value = 0
$ing.Range("G10").NumberFormat = "#,##0.00"
$ing.Range("H10").Value = 0.016
$ing.Range("I10").Value = 0.016
$ing.Range("J10").Value = 2.281
$ing.Range("K10").Value = 0
$ing.Range("K10").NumberFormat = "#,##0.000"
$ing.Range("L10").Formula = "=SUM(F10:K10)"
$ing.Range("L10").NumberFormat = "#,##0.00"
$ing.Range("M10").Formula = "=(C10+E10)/3600"
$ing.Range("N10").Formula = "=(L10+D10)/B10"
$ing.Range("N10").ClearFormats()

# Row 11 - DW10
$ing.Range("A11").Value = "DW10"
$ing.Range("B11").Value = 11794.974
$ing.Range("C11").Value = 14400
$ing.Range("C11").Font.Bold = $true
$ing.Range("C11").Font.Italic = $true
$ing.Range("D11").Value = 152140.85200000001
$ing.Range("D11").NumberFormat = "#,##0.000"
$ing.Range("E11").Value = 10320
$ing.Range("F11").Value = 59608.641000000003
$ing.Range("F11").NumberFormat = "#,##0.00"
$ing.Range("G11").Value = 0
$ing.Range("G11").NumberFormat = "#,##0.00"
$ing.Range("H11").Value = 0.016
$ing.Range("I11").Value = 0.016
$ing.Range("J11").Value = 2.281
$ing.Range("K11").Value = 0
$ing.Range("K11").NumberFormat = "#,##0.000"
$ing.Range("L11").Formula = "=SUM(F11:K11)"
$ing.Range("L11").NumberFormat = "#,##0.00"
$ing.Range("M11").Formula = "=(C11+E11)/3600"
$ing.Range("N11").Formula = "=(L11+D11)/B11"
$ing.Range("N11").ClearFormats()

# Row 12 - DW11
$ing.Range("A12").Value = "DW11"
$ing.Range("B12").Value = 11794.974
$ing.Range("C12").Value = 68400
$ing.Range("C12").Font.Bold = $true
$ing.Range("C12").Font.Italic = $true
$ing.Range("D12").Value = 252155.609
$ing.Range("D12").NumberFormat = "#,##0.000"
$ing.Range("E12").Value = 24015
$ing.Range("E12").NumberFormat = "#,##0.00"
$ing.Range("F12").Value = 49121.366999999998
$ing.Range("F12").NumberFormat = "#,##0.00"
$ing.Range("G12").Value = 0
$ing.Range("G12").NumberFormat = "#,##0.00"
$ing.Range("H12").Value = 0
$ing.Range("H12").NumberFormat = "#,##0.000"
$ing.Range("I12").Value = 0
$ing.Range("I12").NumberFormat = "#,##0.000"
$ing.Range("J12").Value = 0
$ing.Range("J12").NumberFormat = "#,##0.000"
$ing.Range("K12").Value = 0
$ing.Range("K12").NumberFormat = "#,##0.000"
$ing.Range("L12").Formula = "=SUM(F12:K12)"
$ing.Range("L12").NumberFormat = "#,##0.00"
$ing.Range("M12").Formula = "=(C12+E12)/3600"
$ing.Range("M12").ClearFormats()
$ing.Range("N12").Formula = "=(L12+D12)/B12"
$ing.Range("N12").ClearFormats()

# Row 13 - DW12 (partially entered)
$ing.Range("A13").Value = "DW12"
$ing.Range("B13").Value = 11794.974
$ing.Range("C13").Value = 14400
$ing.Range("C13").Font.Bold = $true
$ing.Range("C13").Font.Italic = $true
$ing.Range("D13").Value = 152140.85200000001
$ing.Range("D13").NumberFormat = "#,##0.000"

# Row 14 - just formatted, empty cells (in-progress entry)
$ing.Range("C14").Font.Bold = $true
$ing.Range("C14").Font.Italic = $true
$ing.Range("D14").NumberFormat = "#,##0.000"

# ---------------------------------------------------------------------------
# 2. Ingestion chart: extend series range to cover the new rows
# ---------------------------------------------------------------------------
$chartObj = $ing.ChartObjects(1)
$chart = $chartObj.Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(Ingestion!`$L`$1,Ingestion!`$A`$2:`$A`$12,Ingestion!`$L`$2:`$L`$11,1)"

# ---------------------------------------------------------------------------
# 3. View-state changes recorded for each sheet
# ---------------------------------------------------------------------------

# Performance sheet - selection moved, scrolled back to top
$perf = $wb.Worksheets.Item("Performance")
$perf.Activate()
$perf.Range("B149:O149").Select()

# Elapsed Time sheet - selection moved
$elapsed = $wb.Worksheets.Item("Elapsed Time")
$elapsed.Activate()
$elapsed.Range("J15").Select()

# Total Logical Reads sheet - zoom reset, selection moved, becomes active tab
$reads = $wb.Worksheets.Item("Total Logical Reads")
$reads.Activate()
$excel.ActiveWindow.Zoom = 100
$reads.Range("K18").Select()

# Ingestion sheet - scrolled, selection moved, no longer the active tab
$ing.Activate()
$excel.ActiveWindow.ScrollRow = 4
$ing.Range("I31").Select()

# Re-activate "Total Logical Reads" so it is the workbook's saved active tab
$reads.Activate()

Write-Host "edit applied"
